# Reorder the comma-separated "Recorded By" names in column G for the
# specific rows where the order of entries changed (e.g. the backdoor
# account now listed first, System listed before dnasr281@gmail.com).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$gUpdates = @{
    2 = 'backup@backdoor.com, System, system'
    5 = 'backup@backdoor.com, System'
    8 = 'backup@backdoor.com, System'
    11 = 'System, dnasr281@gmail.com'
    17 = 'System, dnasr281@gmail.com'
    28 = 'backup@backdoor.com, System, system'
    31 = 'backup@backdoor.com, System'
    34 = 'backup@backdoor.com, System'
    37 = 'System, dnasr281@gmail.com'
    43 = 'System, dnasr281@gmail.com'
    54 = 'backup@backdoor.com, System, system'
    57 = 'backup@backdoor.com, System'
    60 = 'backup@backdoor.com, System'
    63 = 'System, dnasr281@gmail.com'
    69 = 'System, dnasr281@gmail.com'
    80 = 'backup@backdoor.com, System'
    81 = 'backup@backdoor.com, System'
    82 = 'backup@backdoor.com, System'
    93 = 'System, dnasr281@gmail.com'
    94 = 'System, dnasr281@gmail.com'
    96 = 'System, dnasr281@gmail.com'
    106 = 'backup@backdoor.com, System'
    107 = 'backup@backdoor.com, System'
    108 = 'backup@backdoor.com, System'
    119 = 'System, dnasr281@gmail.com'
    120 = 'System, dnasr281@gmail.com'
    122 = 'System, dnasr281@gmail.com'
    132 = 'backup@backdoor.com, System'
    133 = 'backup@backdoor.com, System'
    134 = 'backup@backdoor.com, System'
    145 = 'System, dnasr281@gmail.com'
    146 = 'System, dnasr281@gmail.com'
    148 = 'System, dnasr281@gmail.com'
}

foreach ($row in $gUpdates.Keys) {
    $ws.Range("G$row").Value = $gUpdates[$row]
}

